# Data Dictionary.xlsx - 2016.08.18 update
# - Insert a new "Place Level Data?" column (C), shifting old C:G to D:H
# - Mark existing rows as place-level data ("Yes"/"No")
# - Append 5 new reference rows (PM2.5, Smoking, Education, Living Wage x2)
# - Add source-link hyperlinks for the rows that previously lacked one
#   and for all newly appended rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceUrl = "http://www.cdph.ca.gov/programs/Pages/HealthyCommunityIndicators.aspx"

# ---------------------------------------------------------------------
# 1. Insert new column C ("Place Level Data?"), shifting everything right
# ---------------------------------------------------------------------
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "Place Level Data?"
$ws.Columns("C").ColumnWidth = 13.66

# Fill the new column for every existing data row - default "Yes"
# (row 13 is a blank spacer row and must stay empty in column C too)
$placeYesRows = 2..24 | Where-Object { $_ -ne 13 }
foreach ($r in $placeYesRows) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# ---------------------------------------------------------------------
# 2. A12 previously had no source link - give it the standard one
# ---------------------------------------------------------------------
$ws.Range("A12").Value = $sourceUrl
$ws.Hyperlinks.Add($ws.Range("A12"), $sourceUrl) | Out-Null
$ws.Range("A12").Style = $ws.Range("A11").Style

# ---------------------------------------------------------------------
# 3. Append new rows 25-29
# ---------------------------------------------------------------------

# Row 25 - PM2.5 concentration
$ws.Range("A25").Value = $sourceUrl
$ws.Range("B25").Value = "PM25_zcta_place_co_region_ca4-14-13.xls"
$ws.Range("C25").Value = "Yes"
$ws.Range("D25").Value = "Average Ambient PM2.5 Concentration  "
$ws.Range("E25").Value = "2007-2009"
$ws.Range("F25").Value = "No"
$ws.Range("G25").Value = "pm25_conc"
$ws.Range("H25").Value = "Mean ambient PM2.5 concentration (micrograms/cu m.)"

# Row 26 - Smoking prevalence (no place-level data, no year info)
$ws.Range("A26").Value = $sourceUrl
$ws.Range("B26").Value = "HCI_Smoking_755_06NOV15.xlsx"
$ws.Range("C26").Value = "No"
$ws.Range("D26").Value = "Prevalence  of Smoking"
$ws.Range("G26").Value = "N/A"

# Row 27 - Educational attainment
$ws.Range("A27").Value = $sourceUrl
$ws.Range("B27").Value = "ed_attain_ge_hs_output04-14-13.xlsx"
$ws.Range("C27").Value = "Yes"
$ws.Range("D27").Value = "HS or greater educational attainement"
$ws.Range("E27").Value = "2006-2010"
$ws.Range("F27").Value = "2000; 2005-2007; 2008-2010"
$ws.Range("G27").Value = "p_hs_edatt"
$ws.Range("H27").Value = "Percent of population aged 25 years old and over attaining at least a high school degree or GED equivalency"

# Row 28 - Living wage, single mother two children
$ws.Range("A28").Value = $sourceUrl
$ws.Range("B28").Value = "HCI_Living_Wage_770_PL_CO_RE_CA_9-29-13.xls"
$ws.Range("C28").Value = "Yes"
$ws.Range("D28").Value = "Living Wage"
$ws.Range("E28").Value = 2010
$ws.Range("F28").Value = "No"
$ws.Range("G28").Value = "livewage_s"
$ws.Range("H28").Value = "Percent of families living below the living wage - Single mother, two children"

# Row 29 - Living wage, married couple two children
$ws.Range("A29").Value = $sourceUrl
$ws.Range("B29").Value = "HCI_Living_Wage_770_PL_CO_RE_CA_9-29-13.xls"
$ws.Range("C29").Value = "Yes"
$ws.Range("D29").Value = "Living Wage"
$ws.Range("E29").Value = 2010
$ws.Range("F29").Value = "No"
$ws.Range("G29").Value = "livewage_m"
$ws.Range("H29").Value = "Percent of families living below the living wage - Married couple, two children"

# Hyperlinks + styling (copy the look of the existing column-A links) for
# the newly appended rows
foreach ($r in 25..29) {
    $ws.Hyperlinks.Add($ws.Range("A$r"), $sourceUrl) | Out-Null
    $ws.Range("A$r").Style = $ws.Range("A11").Style
}

# ---------------------------------------------------------------------
# 4. Selection / view bookkeeping to match the saved state
# ---------------------------------------------------------------------
$ws.Range("H29").Select()
